$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.378450887226961
$ws.Range("D2").Value = 0.02016070357404942
$ws.Range("E2").Value = 1.48330007287413
$ws.Range("F2").Value = 0.2900348800685251
$ws.Range("G2").Value = 0.1551114414540962
$ws.Range("H2").Value = 0.3165570257724823
$ws.Range("L2").Value = 0.4962858633494989
$ws.Range("M2").Value = 0.3932738228159991
$ws.Range("O2").Value = 0.8429161055172187

$ws.Range("B3").Value = 1.289590310194171
$ws.Range("D3").Value = 0.01757646496097465
$ws.Range("E3").Value = 1.378320081133893
$ws.Range("F3").Value = 0.2880854034570177
$ws.Range("G3").Value = 0.1538616753938484
$ws.Range("H3").Value = 0.3206078687975165
$ws.Range("L3").Value = 0.43635425152587
$ws.Range("M3").Value = 0.356833672571824
$ws.Range("O3").Value = 0.8484505877793111

$ws.Range("B4").Value = 1.23503747975073
$ws.Range("D4").Value = 0.01598351419895749
$ws.Range("E4").Value = 1.314274823987063
$ws.Range("F4").Value = 0.287263024441927
$ws.Range("G4").Value = 0.1533848392322383
$ws.Range("H4").Value = 0.3233820879199882
$ws.Range("L4").Value = 0.3993999589799557
$ws.Range("M4").Value = 0.3343907077836406
$ws.Range("O4").Value = 0.8530553969659138

$ws.Range("B5").Value = 1.212809684819177
$ws.Range("D5").Value = 0.01533287116848925
$ws.Range("E5").Value = 1.288283132790951
$ws.Range("F5").Value = 0.2870215742107831
$ws.Range("G5").Value = 0.1532629311130194
$ws.Range("H5").Value = 0.3245845541237387
$ws.Range("L5").Value = 0.3843023887412755
$ws.Range("M5").Value = 0.3252281936868684
$ws.Range("O5").Value = 0.8552334752999684

$ws.Range("B6").Value = 1.209118981386837
$ws.Range("D6").Value = 0.01522474343553171
$ws.Range("E6").Value = 1.283973814067366
$ws.Range("F6").Value = 0.2869871240148996
$ws.Range("G6").Value = 0.153247042471591
$ws.Range("H6").Value = 0.3247885619513013
$ws.Range("L6").Value = 0.3817931550405831
$ws.Range("M6").Value = 0.3237057617029535
$ws.Range("O6").Value = 0.8556133044441623

$ws.Range("B7").Value = 1.234737694831409
$ws.Range("D7").Value = 0.01597474540256627
$ws.Range("E7").Value = 1.313923851987198
$ws.Range("F7").Value = 0.2872593895646744
$ws.Range("G7").Value = 0.1533829028274809
$ws.Range("H7").Value = 0.3233980137673385
$ws.Range("L7").Value = 0.3991965021918418
$ws.Range("M7").Value = 0.334267206468823
$ws.Range("O7").Value = 0.8530835526392195

$ws.Range("B8").Value = 1.347810637767623
$ws.Range("D8").Value = 0.01927098411536576
$ws.Range("E8").Value = 1.447019286854612
$ws.Range("F8").Value = 0.2892846360923329
$ws.Range("G8").Value = 0.1546198338242348
$ws.Range("H8").Value = 0.3178940636368139
$ws.Range("L8").Value = 0.4756543078970594
$ws.Range("M8").Value = 0.380723724453226
$ws.Range("O8").Value = 0.8445728618186763

$ws.Range("B9").Value = 1.569580963608473
$ws.Range("D9").Value = 0.02568321846760568
$ws.Range("E9").Value = 1.711160991279229
$ws.Range("F9").Value = 0.2962525814909256
$ws.Range("G9").Value = 0.1593800431819119
$ws.Range("H9").Value = 0.3093876028092737
$ws.Range("L9").Value = 0.6243212082545426
$ws.Range("M9").Value = 0.4712675612250763
$ws.Range("O9").Value = 0.837540240735791

$ws.Range("B10").Value = 1.732517640810727
$ws.Range("D10").Value = 0.03036011457328414
$ws.Range("E10").Value = 1.906986666060476
$ws.Range("F10").Value = 0.3032327864335045
$ws.Range("G10").Value = 0.1643415153032208
$ws.Range("H10").Value = 0.3045453061687198
$ws.Range("L10").Value = 0.7327487446042937
$ws.Range("M10").Value = 0.5374401805929097
$ws.Range("O10").Value = 0.8383754767977791

$ws.Range("B11").Value = 1.806639900176378
$ws.Range("D11").Value = 0.03247981036659553
$ws.Range("E11").Value = 1.996425234989886
$ws.Range("F11").Value = 0.3068196539210675
$ws.Range("G11").Value = 0.1669252872059985
$ws.Range("H11").Value = 0.302650868566559
$ws.Range("L11").Value = 0.7818977168363404
$ws.Range("M11").Value = 0.5674663907077644
$ws.Range("O11").Value = 0.8400833280486211

$ws.Range("B12").Value = 1.834707773642378
$ws.Range("D12").Value = 0.03328130265283846
$ws.Range("E12").Value = 2.030341726126323
$ws.Range("F12").Value = 0.3082376508305913
$ws.Range("G12").Value = 0.1679513748997294
$ws.Range("H12").Value = 0.3019780651530226
$ws.Range("L12").Value = 0.8004833982654418
$ws.Range("M12").Value = 0.5788253517045661
$ws.Range("O12").Value = 0.8409229515121694

$ws.Range("B13").Value = 1.828662895343143
$ws.Range("D13").Value = 0.03310874087638638
$ws.Range("E13").Value = 2.023035116380584
$ws.Range("F13").Value = 0.3079295946961622
$ws.Range("G13").Value = 0.1677282576645212
$ws.Range("H13").Value = 0.3021209789513648
$ws.Range("L13").Value = 0.7964818060222854
$ws.Range("M13").Value = 0.57637950521314
$ws.Range("O13").Value = 0.8407335135272547

$ws.Range("B14").Value = 1.808949078077319
$ws.Range("D14").Value = 0.03254577383150092
$ws.Range("E14").Value = 1.999214620339529
$ws.Range("F14").Value = 0.3069351132078921
$ws.Range("G14").Value = 0.1670087448130317
$ws.Range("H14").Value = 0.3025946218313038
$ws.Range("L14").Value = 0.7834272967436391
$ws.Range("M14").Value = 0.5684011281035168
$ws.Range("O14").Value = 0.8401485273130049

$ws.Range("B15").Value = 1.796873691048177
$ws.Range("D15").Value = 0.03220078339528243
$ws.Range("E15").Value = 1.984630041700484
$ws.Range("F15").Value = 0.3063337591759776
$ws.Range("G15").Value = 0.1665742504828955
$ws.Range("H15").Value = 0.3028905543419995
$ws.Range("L15").Value = 0.775427628898683
$ws.Range("M15").Value = 0.5635126575147922
$ws.Range("O15").Value = 0.8398153848879986

$ws.Range("B16").Value = 1.7276735353459
$ws.Range("D16").Value = 0.0302214242328489
$ws.Range("E16").Value = 1.901148541840058
$ws.Range("F16").Value = 0.3030067044501834
$ws.Range("G16").Value = 0.1641792934399717
$ws.Range("H16").Value = 0.3046753373212852
$ws.Range("L16").Value = 0.7295331508894947
$ws.Range("M16").Value = 0.5354763342404425
$ws.Range("O16").Value = 0.8382907520995957

$ws.Range("B17").Value = 1.685221268390592
$ws.Range("D17").Value = 0.02900509713849431
$ws.Range("E17").Value = 1.850024335719837
$ws.Range("F17").Value = 0.3010714609194522
$ws.Range("G17").Value = 0.1627942442006827
$ws.Range("H17").Value = 0.3058493895290297
$ws.Range("L17").Value = 0.7013329328810016
$ws.Range("M17").Value = 0.5182572068165996
$ws.Range("O17").Value = 0.8376969346923886

$ws.Range("B18").Value = 1.660804040804635
$ws.Range("D18").Value = 0.02830476295980588
$ws.Range("E18").Value = 1.820652785764082
$ws.Range("F18").Value = 0.2999970632619551
$ws.Range("G18").Value = 0.1620283175717461
$ws.Range("H18").Value = 0.3065536822193735
$ws.Range("L18").Value = 0.6850964611066956
$ws.Range("M18").Value = 0.5083460838604452
$ws.Range("O18").Value = 0.8374802589713681

$ws.Range("B19").Value = 1.65253684653004
$ws.Range("D19").Value = 0.02806751750431147
$ws.Range("E19").Value = 1.810713978073693
$ws.Range("F19").Value = 0.2996399211878469
$ws.Range("G19").Value = 0.1617742428875175
$ws.Range("H19").Value = 0.3067971186289213
$ws.Range("L19").Value = 0.6795962688813404
$ws.Range("M19").Value = 0.5049891340337211
$ws.Range("O19").Value = 0.837428283740536

$ws.Range("B20").Value = 1.68974037084223
$ws.Range("D20").Value = 0.02913465375608126
$ws.Range("E20").Value = 1.855463126209429
$ws.Range("F20").Value = 0.301273461264195
$ws.Range("G20").Value = 0.1629385010205056
$ws.Range("H20").Value = 0.3057214056172484
$ws.Range("L20").Value = 0.7043366026233286
$ws.Range("M20").Value = 0.5200909536935541
$ws.Range("O20").Value = 0.8377472096089491

$ws.Range("B21").Value = 1.814739527352003
$ws.Range("D21").Value = 0.03271116369941751
$ws.Range("E21").Value = 2.006209998706396
$ws.Range("F21").Value = 0.3072255911660235
$ws.Range("G21").Value = 0.1672187840474919
$ws.Range("H21").Value = 0.302454289676561
$ws.Range("L21").Value = 0.787262433273213
$ws.Range("M21").Value = 0.5707448808248046
$ws.Range("O21").Value = 0.8403151013033323

$ws.Range("B22").Value = 1.896429738566496
$ws.Range("D22").Value = 0.03504165228142142
$ws.Range("E22").Value = 2.10501040992645
$ws.Range("F22").Value = 0.3114640450977291
$ws.Range("G22").Value = 0.1702943627176836
$ws.Range("H22").Value = 0.300578991015243
$ws.Range("L22").Value = 0.8413075637865006
$ws.Range("M22").Value = 0.6037840449406673
$ws.Range("O22").Value = 0.8431186159126298

$ws.Range("B23").Value = 1.852830766309467
$ws.Range("D23").Value = 0.03379848447403333
$ws.Range("E23").Value = 2.052254318947007
$ws.Range("F23").Value = 0.3091698439413904
$ws.Range("G23").Value = 0.1686271926329681
$ws.Range("H23").Value = 0.3015560090910867
$ws.Range("L23").Value = 0.8124767856888866
$ws.Range("M23").Value = 0.5861565905621688
$ws.Range("O23").Value = 0.8415187264621693

$ws.Range("B24").Value = 1.68769731835647
$ws.Range("D24").Value = 0.02907608448064991
$ws.Range("E24").Value = 1.853004185125144
$ws.Range("F24").Value = 0.3011820179608762
$ws.Range("G24").Value = 0.1628731880058609
$ws.Range("H24").Value = 0.3057791758469506
$ws.Range("L24").Value = 0.7029787175899571
$ws.Range("M24").Value = 0.5192619528892806
$ws.Range("O24").Value = 0.8377240919764972

$ws.Range("B25").Value = 1.509584623591024
$ws.Range("D25").Value = 0.02395435060999063
$ws.Range("E25").Value = 1.639383605292693
$ws.Range("F25").Value = 0.2940432337344205
$ws.Range("G25").Value = 0.1578381027993672
$ws.Range("H25").Value = 0.3114426291174155
$ws.Range("L25").Value = 0.5842412662721301
$ws.Range("M25").Value = 0.446833832126714
$ws.Range("O25").Value = 0.8383969894764789
